$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 "I0" and J1 "IF" with the same style as the other headers ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-31: add values for columns I (I0) and J (IF) ---
$data = @(
    @(8,8),
    @(8,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(7,8),
    @(9,9),
    @(3,6),
    @(5,6),
    @(5,7),
    @(6,8),
    @(7,7),
    @(5,7),
    @(1,3),
    @(3,5),
    @(9,9),
    @(5,7),
    @(7,7),
    @(6,8),
    @(8,8),
    @(6,6),
    @(6,6),
    @(5,6),
    @(6,9),
    @(7,7),
    @(4,5),
    @(6,6),
    @(7,7),
    @(1,3),
    @(1,2)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
